# "Basically all done, showed to Richard"
# Mark a few more to-do rows with their status in column C.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5  -> "Instruction note"   => Meh
# Row 11 -> "Axes labels"        => Done
# Row 18 -> "Notes in excel"     => Done
$ws.Range("C5").Value = "Meh"
$ws.Range("C11").Value = "Done"
$ws.Range("C18").Value = "Done"

# Leave the selection where it was left after the last edit.
$ws.Range("C6").Select()
